$wb = $excel.ActiveWorkbook

# --- Sheet 2: rename tab ---------------------------------------------------
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- Sheet 1 ("Metadata") ---------------------------------------------------
$wsMeta = $wb.Worksheets.Item(1)

# URL value: pythia -> cicada
$wsMeta.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/valid-age-reason"

# Date value: updated timestamp
$wsMeta.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# Insert a new row 11 for "Jurisdiction" (pushes Description/Purpose/
# Copyright/Immutable down by one row, as in the target workbook).
$wsMeta.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (style index 2) by
# copying formats down from the row that now sits just below (old row 11,
# now row 12, still carries the original body style).
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMeta.Cells.Item(11, 1).Value = "Jurisdiction"
$wsMeta.Cells.Item(11, 2).Value = ""

# --- Sheet 2 ("Include #0") --------------------------------------------------
$wsInclude.Cells.Item(6, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/ValidAgeReason"
